# Update "想去人数" (interested-count) figures on both the "展览" sheet
# and the aggregated "全部类型" sheet to match newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 139
    $ws.Range("F3").Value = 216
    $ws.Range("F4").Value = 3680
    $ws.Range("F5").Value = 377
}
